$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet / update its tab title to reflect the new date range
$ws.Name = "g3.5b Média(2010-2023)"

# Update region label (row 7 changed from Pernambuco to Amazonas)
$ws.Range("A7").Value = "Amazonas"

# Update the "Valor" column (B) with recalculated averages
$ws.Range("B2").Value = 30.12066375151463
$ws.Range("B3").Value = 16.96368483171459
$ws.Range("B4").Value = 16.38002916428334
$ws.Range("B5").Value = 15.50788422752818
$ws.Range("B6").Value = 14.87449233102686
$ws.Range("B7").Value = 14.71585778266213
$ws.Range("B8").Value = 9.142115924607378
$ws.Range("B9").Value = 14.27646327118632
$ws.Range("B10").Value = 9.852851139694197

# Update the "Ordem" column (C) ranking for row 8
$ws.Range("C8").Value = 19
